$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1 "Save", matching the style of the neighboring header cell G1
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# New data cells H2 and H3 with numeric 0 values
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
